$d = $word.ActiveDocument

# --- Update the date line (keep original run formatting, drop the
#     paragraph mark from the range so we don't split the paragraph) ---
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
$titleRange.Text = "2024-05-25 Saturday"

# --- Update the practice-table numbers ---
# Cell.Range.Text is scoped to the individual cell (unlike Find.Execute on
# a sub-range, which this host treats as a whole-document search), so it
# safely disambiguates the two "433x8=" cells that map to different values.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="209×8="},
    @{Row=1;  Col=2; New="245×2="},
    @{Row=1;  Col=3; New="898×3="},
    @{Row=1;  Col=4; New="421×7="},
    @{Row=1;  Col=5; New="503×4="},

    @{Row=5;  Col=1; New="531×2="},
    @{Row=5;  Col=2; New="806×9="},
    @{Row=5;  Col=3; New="831×2="},
    @{Row=5;  Col=4; New="342×4="},
    @{Row=5;  Col=5; New="770×3="},

    @{Row=10; Col=1; New="950×4="},
    @{Row=10; Col=2; New="870×4="},
    @{Row=10; Col=3; New="950×9="},
    @{Row=10; Col=4; New="865×7="},
    @{Row=10; Col=5; New="147×5="},

    @{Row=15; Col=1; New="402×5="},
    @{Row=15; Col=2; New="861×2="},
    @{Row=15; Col=3; New="430×8="},
    @{Row=15; Col=4; New="726×2="},
    @{Row=15; Col=5; New="243×3="},

    @{Row=20; Col=1; New="210×5="},
    @{Row=20; Col=2; New="868×6="},
    @{Row=20; Col=3; New="996×8="},
    @{Row=20; Col=4; New="225×8="},
    @{Row=20; Col=5; New="604×7="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    $cellEnd = $cellRange.End - 1
    $valueRange = $d.Range($cellRange.Start, $cellEnd)
    $valueRange.Text = $u.New
}
